# Append two new data rows (115, 116) to the bottom of the "particelle non
# trovate" table on Sheet1, matching the pattern of the existing rows:
#   col A -> sequential index, styled like the rest of column A
#   col B -> codice_particella (plain text, even when it looks numeric)
#   col C -> codice_comune_catastale (plain number, unstyled)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, [string]$text) {
    # Force a literal/text value even for strings that look numeric
    # (e.g. ".451", "1768"), then drop the number-format style that
    # gets stamped on so the cell ends up un-styled, like the source rows.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 115: index 113, particella ".451", comune 88
$ws.Cells.Item(115, 1).Value = 113
$ws.Range("A114").Copy()
$ws.Range("A115").PasteSpecial(-4122)  # xlPasteFormats: copy A114's style (index "1")
Set-TextCell $ws.Cells.Item(115, 2) ".451"
$ws.Cells.Item(115, 3).Value = 88

# Row 116: index 114, particella "1768", comune 240
$ws.Cells.Item(116, 1).Value = 114
$ws.Range("A114").Copy()
$ws.Range("A116").PasteSpecial(-4122)
Set-TextCell $ws.Cells.Item(116, 2) "1768"
$ws.Cells.Item(116, 3).Value = 240

$excel.CutCopyMode = 0
